# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> bound to the (one) slide master, i.e. the theme
#                            actually applied across the slides ("Integral" /
#                            "Red Violet" color scheme before the edit).
#   ppt/theme/theme2.xml -> bound to the notes master ("Office Theme" /
#                            "Office" color scheme before the edit).
#
# The commit swaps the full contents of the two parts: the slide master's
# theme becomes the plain "Office Theme" color scheme, and the notes
# master's theme becomes the former "Integral" color scheme. The font
# scheme and format scheme (fills/lines/effects) are identical between the
# two theme parts, so the only real content difference is the color scheme
# (clrScheme) -- re-pointing the 12 theme colors reproduces the swap for the
# part of the theme that actually drives the presentation's look.

$p = $ppt.ActivePresentation

$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
